$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Remove the "Meta description" paragraph that currently sits right
#        after the title heading (empty run + bold "Meta description" run +
#        plain run with the description text).
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        [void]$p.Range.Delete()
        break
    }
}

# --- 2. Insert a new paragraph, with bold text repeating the page title,
#        right before the final (italic, AI image-prompt) paragraph.
$last = $d.Paragraphs.Last
[void]$last.Range.InsertParagraphBefore()

$cnt = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($cnt - 1)
[void]$newPara.Range.InsertXML(
  "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr>" +
  "<w:t>Play Dim Sum Prize for Free - Delicious Chinese-themed Slot</w:t>" +
  "</w:r></w:p>"
)

# --- 3. Replace the final paragraph's text (the old AI image prompt) with
#        the meta-description copy, keeping its italic formatting.
$last = $d.Paragraphs.Last
[void]$last.Range.InsertXML(
  "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr>" +
  "<w:t>Dim Sum Prize is a Chinese-themed slot game with 10 fixed paylines. Play now for free and enjoy two bonus features and an impressive RTP of 97.18%.</w:t>" +
  "</w:r></w:p>"
)
